# Commit: "#5: property aircraft done"
#
# The "航空器" (aircraft) sheet only ever held bogus placeholder data
# (currency-disclosure headers such as 幣別/所有人/外幣總額/新臺幣總額或折合新臺幣總額
# with no real row), so it is removed entirely now that aircraft property
# has been confirmed to not apply. Removing it frees up those four now-unused
# shared strings.
#
# Separately, the "建物" (building) sheet's property_category column (I)
# was mistakenly tagged with the "land" category string; fix it to read
# "building" for every data row.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Fix the property_category column on the 建物 (building) sheet: every
# data row (2-22) incorrectly says "land" instead of "building".
$building = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 22; $r++) {
    $building.Cells.Item($r, 9).Value = "building"
}

# Drop the placeholder 航空器 (aircraft) sheet entirely.
$aircraft = $wb.Worksheets.Item("航空器")
$aircraft.Delete()
